# Implement CXR pre-Xpert alternative for ACF
#
# Adds a new parameter row "tb_sensitivity_cxr" (value 0.9) just above the
# existing "program_timeperiod_await_treatment_smearneg_xpert" row (row 37)
# on the "constants" worksheet. All rows from the old row 37 onward shift
# down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Insert a new blank row at position 37; this pushes all subsequent rows
# down by one and inherits the formatting of the row above (row 36), which
# matches the desired style for the new label/value pair.
$ws.Rows.Item(37).Insert()

# Populate the new row: parameter name in column A, point value in column B.
$ws.Cells.Item(37, 1).Value = "tb_sensitivity_cxr"
$ws.Cells.Item(37, 2).Value = 0.9

# Update the active selection to match the author's final cursor position.
$ws.Activate()
$ws.Range("B39").Select()
